$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.877.57'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.311.77'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.21'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.68'
$ws.Range("E6").Value = '  +1.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.33'
$ws.Range("E10").Value = '  +3.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.63'
$ws.Range("E12").Value = '  +4.86%  '
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.982'
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.46'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '2.661.96'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").Value = '2.310.41'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '42.815.55'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.57'
$ws.Range("E21").Value = '  +33.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.05'
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.56'
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '271.58'
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.02'
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("E28").Value = '  -3.34%  '
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.31'
$ws.Range("E30").Value = '  +8.85%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.24'
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.25'
$ws.Range("E32").Value = '  +7.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0893'
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.56'
$ws.Range("E35").Value = '  -10.28%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.116'
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.64'
$ws.Range("E37").Value = '  +2.21%  '
$ws.Range("E38").Value = '  +3.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.74'
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("E40").Value = '  -4.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.60'
$ws.Range("E41").Value = '  +10.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.78'
$ws.Range("E42").Value = '  +1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.96'
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.42'
$ws.Range("E46").Value = '  +4.98%  '
$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '82.94'
$ws.Range("E47").Value = '  +9.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '115.30'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.32'
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").Value = '1.609.00'
$ws.Range("E51").Value = '  +5.46%  '
